# df_list_20241219.xlsx — "Add files via upload" re-upload edit.
#
# The commit re-uploads the workbook after the author scrolled/reselected
# and re-saved roughly 12 hours later. The only semantically meaningful
# content change is that the "collection timestamp" column (F) for every
# data row (rows 2-77) moved forward by exactly half a day (12h) — e.g.
# 2024-12-19 00:25:00 -> 2024-12-19 12:25:00 — while everything else
# (labels, URLs, the E-column dates, styles) stayed the same. The sheet's
# on-screen selection also moved from L78 to H76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the F-column "collected at" timestamp for every data row (2-77)
# forward by 0.5 days (12 hours), preserving the exact fractional value
# already stored in each cell.
for ($row = 2; $row -le 77; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $current = $cell.Value2()
    $cell.Value2 = $current + 0.5
}

# Match the author's final on-screen selection / scroll position.
$ws.Range("H76").Select() | Out-Null
